# Weekly update: insert the newest Níspero (Vega Modelo de Temuco) price row
# above the existing historical rows, pushing the previous rows (10-14) down
# to (11-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10; existing rows 10-14 shift down to 11-15.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's data.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44874
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100104
$ws.Range("H10").Value = "Frutos de pepita"
$ws.Range("I10").Value = 100104004
$ws.Range("J10").Value = "Níspero"
$ws.Range("K10").Value = "Californiana(o)"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 25000
$ws.Range("O10").Value = 25000
$ws.Range("P10").Value = 25000
$ws.Range("Q10").Value = "`$/bandeja 10 kilos"
$ws.Range("R10").Value = "Provincia de Quillota"
$ws.Range("S10").Value = 2500
$ws.Range("T10").Value = 10
